$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 8; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $textA = $cellA.Text
    if ($textA -match '220120') {
        $cellA.Value = $textA -replace '220120', '230120'
    }

    $cellB = $ws.Cells.Item($r, 2)
    $textB = $cellB.Text
    if ($textB -match '220120') {
        $cellB.Value = $textB -replace '220120', '230120'
    }
}
